# Updated cryptos list on Tue Nov 14 22:33:01 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "35.584.67"
$ws.Range("E2").Value = "  -2.76%  "

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "1.982.00"
$ws.Range("E3").Value = "  -4.09%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  +0.02%  "

# --- Row 5 (BNB) ---
$ws.Range("D5").Value = "'241.48"
$ws.Range("E5").Value = "  -0.44%  "

# --- Row 6 (XRP) ---
$ws.Range("D6").Value = "'0.636"
$ws.Range("E6").Value = "  -3.67%  "

# --- Row 8 (Solana) ---
$ws.Range("D8").Value = "'56.21"
$ws.Range("E8").Value = "  +6.52%  "

# --- Row 9 (OKB) ---
$ws.Range("D9").Value = "'60.16"
$ws.Range("E9").Value = "  +1.81%  "

# --- Row 10 (Cardano) ---
$ws.Range("E10").Value = "  -0.43%  "

# --- Row 11 (Dogecoin) ---
$ws.Range("E11").Value = "  -3.17%  "

# --- Row 12 (TRON) ---
$ws.Range("E12").Value = "  -5.14%  "

# --- Row 13 (Polygon) ---
$ws.Range("D13").Value = "'0.910"
$ws.Range("E13").Value = "  +0.47%  "

# --- Row 14 ---
$ws.Range("D14").Value = "'14.13"
$ws.Range("E14").Value = "  -3.67%  "

# --- Row 15 ---
$ws.Range("D15").Value = "2.273.83"
$ws.Range("E15").Value = "  -3.90%  "

# --- Row 16 ---
$ws.Range("E16").Value = "  -3.36%  "

# --- Row 17 ---
$ws.Range("D17").Value = "1.981.33"
$ws.Range("E17").Value = "  -4.39%  "

# --- Row 18 ---
$ws.Range("D18").Value = "'16.99"
$ws.Range("E18").Value = "  +3.68%  "

# --- Row 19 ---
$ws.Range("D19").Value = "35.491.93"
$ws.Range("E19").Value = "  -2.79%  "

# --- Row 20 ---
$ws.Range("D20").Value = "'70.17"
$ws.Range("E20").Value = "  -2.24%  "

# --- Row 21 ---
$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  -3.25%  "

# --- Row 22 ---
$ws.Range("D22").Value = "'232.17"
$ws.Range("E22").Value = "  -2.15%  "

# --- Row 23 ---
$ws.Range("E23").Value = "  -4.18%  "

# --- Row 24 ---
$ws.Range("E24").Value = "  -0.15%  "

# --- Row 25 ---
$ws.Range("E25").Value = "  -3.14%  "

# --- Row 26 ---
$ws.Range("D26").Value = "'2.31"
$ws.Range("E26").Value = "  +8.45%  "

# --- Row 27 ---
$ws.Range("D27").Value = "'163.43"
$ws.Range("E27").Value = "  -0.60%  "

# --- Row 28 ---
$ws.Range("E28").Value = "  -4.32%  "

# --- Row 29 ---
$ws.Range("D29").Value = "'19.42"
$ws.Range("E29").Value = "  -5.58%  "

# --- Row 30 ---
$ws.Range("D30").Value = "'0.119"
$ws.Range("E30").Value = "  -2.81%  "

# --- Row 31 ---
$ws.Range("E31").Value = "  -5.44%  "

# --- Row 32 ---
$ws.Range("E32").Value = "  -1.89%  "

# --- Row 33 ---
$ws.Range("D33").Value = "'0.0584"
$ws.Range("E33").Value = "  -2.13%  "

# --- Row 34 ---
$ws.Range("D34").Value = "'0.0908"
$ws.Range("E34").Value = "  +11.08%  "

# --- Row 35 ---
$ws.Range("D35").Value = "'4.24"
$ws.Range("E35").Value = "  -7.29%  "

# --- Row 36: BinanceUSD -> LidoDAOToken ---
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.32"
$ws.Range("E36").Value = "  +2.15%  "

# --- Row 37: LidoDAOToken -> BinanceUSD ---
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.04%  "

# --- Row 38 ---
$ws.Range("E38").Value = "  -2.25%  "

# --- Row 39 ---
$ws.Range("D39").Value = "'4.88"
$ws.Range("E39").Value = "  +0.33%  "

# --- Row 40 ---
$ws.Range("E40").Value = "  -5.21%  "

# --- Row 41 ---
$ws.Range("D41").Value = "'2.82"
$ws.Range("E41").Value = "  -3.61%  "

# --- Row 42 ---
$ws.Range("E42").Value = "  -3.18%  "

# --- Row 43 ---
$ws.Range("E43").Value = "  -4.88%  "

# --- Row 44 ---
$ws.Range("D44").Value = "'0.0885"
$ws.Range("E44").Value = "  -5.16%  "

# --- Row 45 ---
$ws.Range("D45").Value = "'90.48"
$ws.Range("E45").Value = "  -4.11%  "

# --- Row 46 ---
$ws.Range("D46").Value = "1.373.43"
$ws.Range("E46").Value = "  -1.40%  "

# --- Row 47 ---
$ws.Range("D47").Value = "'7.32"
$ws.Range("E47").Value = "  -1.30%  "

# --- Row 48 ---
$ws.Range("E48").Value = "  -0.83%  "

# --- Row 49 ---
$ws.Range("E49").Value = "  +0.68%  "

# --- Row 50 ---
$ws.Range("E50").Value = "  -4.01%  "

# --- Row 51 ---
$ws.Range("D51").Value = "'45.59"
$ws.Range("E51").Value = "  +0.45%  "
